$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Note: cells that hold text which *looks* numeric/date/percentage-like
# (e.g. "69.6%", "2026-02-16", "3.69") are written with NumberFormat "@"
# (Text) set first, so Excel's auto-conversion doesn't turn them into a real
# number/date. Plain words/sentences and genuine numeric values are written
# directly.
# ---------------------------------------------------------------------------

# ===========================================================================
# Summary sheet
# ===========================================================================
$summary = $wb.Worksheets.Item("Summary")

# Row 2 - OVERALL / ALL COMBINED
$summary.Cells.Item(2,3).Value = 23
$summary.Cells.Item(2,4).NumberFormat = "@"
$summary.Cells.Item(2,4).Value = "69.6%"
$summary.Cells.Item(2,5).NumberFormat = "@"
$summary.Cells.Item(2,5).Value = "+5.2943%"
$summary.Cells.Item(2,6).NumberFormat = "@"
$summary.Cells.Item(2,6).Value = "+0.2302%"

# Row 4 - STRATEGY / momentum
$summary.Cells.Item(4,3).Value = 11
$summary.Cells.Item(4,4).NumberFormat = "@"
$summary.Cells.Item(4,4).Value = "36.4%"
$summary.Cells.Item(4,5).NumberFormat = "@"
$summary.Cells.Item(4,5).Value = "+1.7430%"
$summary.Cells.Item(4,6).NumberFormat = "@"
$summary.Cells.Item(4,6).Value = "+0.1585%"

# ===========================================================================
# leadlag sheet - append trade #51 (new row 41), still OPEN
# ===========================================================================
$leadlag = $wb.Worksheets.Item("leadlag")

$leadlag.Cells.Item(41,1).Value = 51
$leadlag.Cells.Item(41,2).NumberFormat = "@"
$leadlag.Cells.Item(41,2).Value = "2026-02-16"
$leadlag.Cells.Item(41,3).Value = "21:31:02"
$leadlag.Cells.Item(41,4).Value = "leadlag"
$leadlag.Cells.Item(41,5).Value = "DOWN"
$leadlag.Cells.Item(41,6).Value = 68754.095
# column G (Exit Price) stays blank - trade is still OPEN
$leadlag.Cells.Item(41,8).Value = "OPEN"
$leadlag.Cells.Item(41,9).Value = 0
$leadlag.Cells.Item(41,10).Value = 0
$leadlag.Cells.Item(41,11).Value = 0.7212
$leadlag.Cells.Item(41,12).Value = "Coinbase leading with -0.072% move"
# column M (Exit Reason) stays blank - trade is still OPEN
$leadlag.Cells.Item(41,14).Value = 0

# ===========================================================================
# momentum sheet - close trades #21, #22, #23 (rows 4, 5, 6)
# ===========================================================================
$momentum = $wb.Worksheets.Item("momentum")

$momentum.Cells.Item(4,7).Value = 68726.060364
$momentum.Cells.Item(4,8).Value = "CLOSED"
$momentum.Cells.Item(4,9).Value = 0.5087
$momentum.Cells.Item(4,10).Value = 5.09
$momentum.Cells.Item(4,13).Value = "time_exit_5min"
$momentum.Cells.Item(4,14).Value = 5

$momentum.Cells.Item(5,7).Value = 68454.970218
$momentum.Cells.Item(5,8).Value = "CLOSED"
$momentum.Cells.Item(5,9).Value = 0.9359
$momentum.Cells.Item(5,10).Value = 9.359999999999999
$momentum.Cells.Item(5,13).Value = "time_exit_5min"
$momentum.Cells.Item(5,14).Value = 5

$momentum.Cells.Item(6,7).Value = 68916.001487
$momentum.Cells.Item(6,8).Value = "CLOSED"
$momentum.Cells.Item(6,9).Value = 0.2531
$momentum.Cells.Item(6,10).Value = 2.53
$momentum.Cells.Item(6,13).Value = "time_exit_5min"
$momentum.Cells.Item(6,14).Value = 5

# ===========================================================================
# All Trades sheet - append the 3 newly-closed momentum trades as rows 22-24
# ===========================================================================
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Cells.Item(22,1).Value = 21
$allTrades.Cells.Item(22,2).NumberFormat = "@"
$allTrades.Cells.Item(22,2).Value = "2026-02-16"
$allTrades.Cells.Item(22,3).Value = "21:25:35"
$allTrades.Cells.Item(22,4).Value = "momentum"
$allTrades.Cells.Item(22,5).Value = "DOWN"
$allTrades.Cells.Item(22,6).Value = 69077.44500000001
$allTrades.Cells.Item(22,7).Value = 68726.060364
$allTrades.Cells.Item(22,8).Value = "CLOSED"
$allTrades.Cells.Item(22,9).Value = 0.5087
$allTrades.Cells.Item(22,10).Value = 5.09
$allTrades.Cells.Item(22,11).Value = 0.9
$allTrades.Cells.Item(22,12).Value = "Downward momentum: -0.283% over 10 samples"
$allTrades.Cells.Item(22,13).Value = "time_exit_5min"
$allTrades.Cells.Item(22,14).Value = 5

$allTrades.Cells.Item(23,1).Value = 22
$allTrades.Cells.Item(23,2).NumberFormat = "@"
$allTrades.Cells.Item(23,2).Value = "2026-02-16"
$allTrades.Cells.Item(23,3).Value = "21:25:41"
$allTrades.Cells.Item(23,4).Value = "momentum"
$allTrades.Cells.Item(23,5).Value = "DOWN"
$allTrades.Cells.Item(23,6).Value = 69101.72500000001
$allTrades.Cells.Item(23,7).Value = 68454.970218
$allTrades.Cells.Item(23,8).Value = "CLOSED"
$allTrades.Cells.Item(23,9).Value = 0.9359
$allTrades.Cells.Item(23,10).Value = 9.359999999999999
$allTrades.Cells.Item(23,11).Value = 0.9
$allTrades.Cells.Item(23,12).Value = "Downward momentum: -0.262% over 10 samples"
$allTrades.Cells.Item(23,13).Value = "time_exit_5min"
$allTrades.Cells.Item(23,14).Value = 5

$allTrades.Cells.Item(24,1).Value = 23
$allTrades.Cells.Item(24,2).NumberFormat = "@"
$allTrades.Cells.Item(24,2).Value = "2026-02-16"
$allTrades.Cells.Item(24,3).Value = "21:25:47"
$allTrades.Cells.Item(24,4).Value = "momentum"
$allTrades.Cells.Item(24,5).Value = "DOWN"
$allTrades.Cells.Item(24,6).Value = 69090.855
$allTrades.Cells.Item(24,7).Value = 68916.001487
$allTrades.Cells.Item(24,8).Value = "CLOSED"
$allTrades.Cells.Item(24,9).Value = 0.2531
$allTrades.Cells.Item(24,10).Value = 2.53
$allTrades.Cells.Item(24,11).Value = 0.9
$allTrades.Cells.Item(24,12).Value = "Downward momentum: -0.280% over 10 samples"
$allTrades.Cells.Item(24,13).Value = "time_exit_5min"
$allTrades.Cells.Item(24,14).Value = 5

# ===========================================================================
# Comparison sheet - row 3 (momentum) stats
# ===========================================================================
$comparison = $wb.Worksheets.Item("Comparison")

$comparison.Cells.Item(3,2).Value = 11
$comparison.Cells.Item(3,3).NumberFormat = "@"
$comparison.Cells.Item(3,3).Value = "36.4%"
$comparison.Cells.Item(3,4).NumberFormat = "@"
$comparison.Cells.Item(3,4).Value = "3.69"
$comparison.Cells.Item(3,5).NumberFormat = "@"
$comparison.Cells.Item(3,5).Value = "+0.5976%"
$comparison.Cells.Item(3,7).NumberFormat = "@"
$comparison.Cells.Item(3,7).Value = "0.92"
